$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (A3) switches from a date-only format to a date+time format,
# matching the format already used on A2.
$ws.Range("A3").NumberFormat = $ws.Range("A2").NumberFormat

# New row 4: daily update values, with A4 taking on the date-only format
# that A3 previously had.
$ws.Range("A4").Value = 45953
$ws.Range("B4").Value = 4
$ws.Range("C4").Value = 8
$ws.Range("D4").Value = 4
$ws.Range("A4").NumberFormat = "YYYY-MM-DD"
